$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.352688459262481
$ws.Range("C2").Value = 0.2438158961542172
$ws.Range("D2").Value = 0.5228266803324146
$ws.Range("E2").Value = 0.1829047581761003
$ws.Range("G2").Value = 1.032520994063276
$ws.Range("H2").Value = 1.057146006784222
$ws.Range("J2").Value = 0.08816575100385293
$ws.Range("L2").Value = 0.3668636860179362
$ws.Range("M2").Value = 0.3458638963896732
$ws.Range("O2").Value = 4.237916805190395
# Row 3
$ws.Range("B3").Value = 1.252488878539054
$ws.Range("C3").Value = 0.2332425560761209
$ws.Range("D3").Value = 0.524071164148836
$ws.Range("E3").Value = 0.184498590239718
$ws.Range("G3").Value = 1.040758594272347
$ws.Range("H3").Value = 1.06613384883002
$ws.Range("J3").Value = 0.08781546449406008
$ws.Range("L3").Value = 0.3627348923325329
$ws.Range("M3").Value = 0.329157898387308
$ws.Range("O3").Value = 4.273857945568153
# Row 4
$ws.Range("B4").Value = 1.191100278037538
$ws.Range("C4").Value = 0.2267397842876591
$ws.Range("D4").Value = 0.5250621559068662
$ws.Range("E4").Value = 0.1855388112822123
$ws.Range("G4").Value = 1.046515457162741
$ws.Range("H4").Value = 1.072152772717168
$ws.Range("J4").Value = 0.0876074842823229
$ws.Range("L4").Value = 0.3603191279266298
$ws.Range("M4").Value = 0.3189734253100625
$ws.Range("O4").Value = 4.298443277077595
# Row 5
$ws.Range("B5").Value = 1.166119396706051
$ws.Range("C5").Value = 0.2240873777530226
$ws.Range("D5").Value = 0.5255231236420173
$ws.Range("E5").Value = 0.1859782185124965
$ws.Range("G5").Value = 1.049037034540767
$ws.Range("H5").Value = 1.074731393870479
$ws.Range("J5").Value = 0.0875245328988612
$ws.Range("L5").Value = 0.3593648214612628
$ws.Range("M5").Value = 0.3148418501362329
$ws.Range("O5").Value = 4.309094774248919
# Row 6
$ws.Range("B6").Value = 1.16197353278119
$ws.Range("C6").Value = 0.2236468046298228
$ws.Range("D6").Value = 0.5256031199210867
$ws.Range("E6").Value = 0.1860521189788464
$ws.Range("G6").Value = 1.049466342612234
$ws.Range("H6").Value = 1.075167174902106
$ws.Range("J6").Value = 0.0875108682566399
$ws.Range("L6").Value = 0.359208183653692
$ws.Range("M6").Value = 0.3141569424729056
$ws.Range("O6").Value = 4.310901658634975
# Row 7
$ws.Range("B7").Value = 1.190763231264157
$ws.Range("C7").Value = 0.226704022736044
$ws.Range("D7").Value = 0.525068141250685
$ws.Range("E7").Value = 0.1855446744630909
$ws.Range("G7").Value = 1.046548753198515
$ws.Range("H7").Value = 1.072187039242998
$ws.Range("J7").Value = 0.08760635825001373
$ws.Range("L7").Value = 0.3603061356184298
$ws.Range("M7").Value = 0.3189176294013549
$ws.Range("O7").Value = 4.298584365120817
# Row 8
$ws.Range("B8").Value = 1.318112737704212
$ws.Range("C8").Value = 0.2401725544815747
$ws.Range("D8").Value = 0.523208738660955
$ws.Range("E8").Value = 0.1834415372423939
$ws.Range("G8").Value = 1.035216164606965
$ws.Range("H8").Value = 1.060141192430336
$ws.Range("J8").Value = 0.08804350957399976
$ws.Range("L8").Value = 0.365415386700235
$ws.Range("M8").Value = 0.3400886775295646
$ws.Range("O8").Value = 4.249786701366844
# Row 9
$ws.Range("B9").Value = 1.568849586172973
$ws.Range("C9").Value = 0.2664912944557898
$ws.Range("D9").Value = 0.5213594073881467
$ws.Range("E9").Value = 0.1798052137390833
$ws.Range("G9").Value = 1.01854590073188
$ws.Range("H9").Value = 1.040487615229083
$ws.Range("J9").Value = 0.08895636567446985
$ws.Range("L9").Value = 0.3763765405256976
$ws.Range("M9").Value = 0.3821737939136085
$ws.Range("O9").Value = 4.17408018923814
# Row 10
$ws.Range("B10").Value = 1.753615635340793
$ws.Range("C10").Value = 0.285762055144346
$ws.Range("D10").Value = 0.5210924540385378
$ws.Range("E10").Value = 0.177429829473045
$ws.Range("G10").Value = 1.009694063873212
$ws.Range("H10").Value = 1.028465085958928
$ws.Range("J10").Value = 0.08966008125729985
$ws.Range("L10").Value = 0.3849985087501011
$ws.Range("M10").Value = 0.4134289469066061
$ws.Range("O10").Value = 4.130660683796833
# Row 11
$ws.Range("B11").Value = 1.837777888375683
$ws.Range("C11").Value = 0.2945128426402732
$ws.Range("D11").Value = 0.5212073263998604
$ws.Range("E11").Value = 0.1764132657848041
$ws.Range("G11").Value = 1.006406830602089
$ws.Range("H11").Value = 1.02352010522624
$ws.Range("J11").Value = 0.08998722077528143
$ws.Range("L11").Value = 0.3890433133595934
$ws.Range("M11").Value = 0.4277183130255295
$ws.Range("O11").Value = 4.113561716722927
# Row 12
$ws.Range("B12").Value = 1.86966253583887
$ws.Range("C12").Value = 0.2978241141850049
$ws.Range("D12").Value = 0.5212847354580106
$ws.Range("E12").Value = 0.1760375054613403
$ws.Range("G12").Value = 1.005268560685053
$ws.Range("H12").Value = 1.021722907794171
$ws.Range("J12").Value = 0.09011209284922472
$ws.Range("L12").Value = 0.3905924970174937
$ws.Range("M12").Value = 0.4331393135741308
$ws.Range("O12").Value = 4.107468572696945
# Row 13
$ws.Range("B13").Value = 1.862795001706672
$ws.Range("C13").Value = 0.2971110856023813
$ws.Range("D13").Value = 0.5212665570465589
$ws.Range("E13").Value = 0.176118023657553
$ws.Range("G13").Value = 1.00550896628863
$ws.Range("H13").Value = 1.022106614934444
$ws.Range("J13").Value = 0.09008515559493091
$ws.Range("L13").Value = 0.3902580761656367
$ws.Range("M13").Value = 0.4319713675902221
$ws.Range("O13").Value = 4.108763850498804
# Row 14
$ws.Range("B14").Value = 1.840400782428844
$ws.Range("C14").Value = 0.2947853135417233
$ws.Range("D14").Value = 0.5212130157068344
$ws.Range("E14").Value = 0.176382167651802
$ws.Range("G14").Value = 1.006311048249444
$ws.Range("H14").Value = 1.023370738449941
$ws.Range("J14").Value = 0.08999747430492278
$ws.Range("L14").Value = 0.3891704155676763
$ws.Range("M14").Value = 0.4281641048421463
$ws.Range("O14").Value = 4.113052774715442
# Row 15
$ws.Range("B15").Value = 1.8266854772491
$ws.Range("C15").Value = 0.2933603838707199
$ws.Range("D15").Value = 0.5211846339626192
$ws.Range("E15").Value = 0.1765451598969943
$ws.Range("G15").Value = 1.006816226189002
$ws.Range("H15").Value = 1.024154864349768
$ws.Range("J15").Value = 0.08994389558374039
$ws.Range("L15").Value = 0.3885064677339614
$ws.Range("M15").Value = 0.4258333300023693
$ws.Range("O15").Value = 4.115729606299169
# Row 16
$ws.Range("B16").Value = 1.748117521292386
$ws.Range("C16").Value = 0.2851898381233013
$ws.Range("D16").Value = 0.5210896958964497
$ws.Range("E16").Value = 0.177497551127126
$ws.Range("G16").Value = 1.009923789699272
$ws.Range("H16").Value = 1.028798798056812
$ws.Range("J16").Value = 0.08963884160202795
$ws.Range("L16").Value = 0.3847366277315984
$ws.Range("M16").Value = 0.4124965096509428
$ws.Range("O16").Value = 4.13183155625407
# Row 17
$ws.Range("B17").Value = 1.699945856048089
$ws.Range("C17").Value = 0.2801733244908746
$ws.Range("D17").Value = 0.5210919282051378
$ws.Range("E17").Value = 0.1780981966869035
$ws.Range("G17").Value = 1.012019723733616
$ws.Range("H17").Value = 1.031781935904377
$ws.Range("J17").Value = 0.08945348531855046
$ws.Range("L17").Value = 0.3824552683750255
$ws.Range("M17").Value = 0.4043328227409759
$ws.Range("O17").Value = 4.142389305483249
# Row 18
$ws.Range("B18").Value = 1.672249364570462
$ws.Range("C18").Value = 0.2772865053229339
$ws.Range("D18").Value = 0.5211154591499252
$ws.Range("E18").Value = 0.1784496979521686
$ws.Range("G18").Value = 1.013294844426525
$ws.Range("H18").Value = 1.033547091863184
$ws.Range("J18").Value = 0.0893475346735606
$ws.Range("L18").Value = 0.3811546408599042
$ws.Range("M18").Value = 0.3996440102005892
$ws.Range("O18").Value = 4.148711512561817
# Row 19
$ws.Range("B19").Value = 1.662873679060453
$ws.Range("C19").Value = 0.2763088351612168
$ws.Range("D19").Value = 0.5211272494528174
$ws.Range("E19").Value = 0.178569745621413
$ws.Range("G19").Value = 1.013738525468952
$ws.Range("H19").Value = 1.034153216795033
$ws.Range("J19").Value = 0.0893117757336519
$ws.Range("L19").Value = 0.3807162584294161
$ws.Range("M19").Value = 0.3980576234176922
$ws.Range("O19").Value = 4.150894972024275
# Row 20
$ws.Range("B20").Value = 1.705072730301481
$ws.Range("C20").Value = 0.2807074929462487
$ws.Range("D20").Value = 0.5210893886143566
$ws.Range("E20").Value = 0.1780336333955468
$ws.Range("G20").Value = 1.011789403765121
$ws.Range("H20").Value = 1.031459269974874
$ws.Range("J20").Value = 0.08947314848706611
$ws.Range("L20").Value = 0.3826969285295263
$ws.Range("M20").Value = 0.4052011675148819
$ws.Range("O20").Value = 4.141239572365521
# Row 21
$ws.Range("B21").Value = 1.846978135337963
$ws.Range("C21").Value = 0.2954685177396073
$ws.Range("D21").Value = 0.5212278223767299
$ws.Range("E21").Value = 0.1763043328499441
$ws.Range("G21").Value = 1.006072564500229
$ws.Range("H21").Value = 1.022997389571628
$ws.Range("J21").Value = 0.09002320165559397
$ws.Range("L21").Value = 0.3894894138301055
$ws.Range("M21").Value = 0.4292821232768063
$ws.Range("O21").Value = 4.111782646885871
# Row 22
$ws.Range("B22").Value = 1.939803646004236
$ws.Range("C22").Value = 0.3051012506558948
$ws.Range("D22").Value = 0.5215159134337739
$ws.Range("E22").Value = 0.1752276963880366
$ws.Range("G22").Value = 1.002957329886669
$ws.Range("H22").Value = 1.017906309863505
$ws.Range("J22").Value = 0.0903884641066206
$ws.Range("L22").Value = 0.3940306622194782
$ws.Range("M22").Value = 0.4450781409718374
$ws.Range("O22").Value = 4.094756819736091
# Row 23
$ws.Range("B23").Value = 1.890253990409633
$ws.Range("C23").Value = 0.2999614745662029
$ws.Range("D23").Value = 0.5213440955903224
$ws.Range("E23").Value = 0.1757974213739537
$ws.Range("G23").Value = 1.004563099597391
$ws.Range("H23").Value = 1.020583326762036
$ws.Range("J23").Value = 0.09019299433256123
$ws.Range("L23").Value = 0.3915976252497018
$ws.Range("M23").Value = 0.4366423315770334
$ws.Range("O23").Value = 4.103640016624752
# Row 24
$ws.Range("B24").Value = 1.70275487642175
$ws.Range("C24").Value = 0.280466003962772
$ws.Range("D24").Value = 0.5210904674593877
$ws.Range("E24").Value = 0.1780628032016018
$ws.Range("G24").Value = 1.011893312953148
$ws.Range("H24").Value = 1.031604991147958
$ws.Range("J24").Value = 0.08946425685737935
$ws.Range("L24").Value = 0.3825876398369559
$ws.Range("M24").Value = 0.4048085744911916
$ws.Range("O24").Value = 4.141758580467865
# Row 25
$ws.Range("B25").Value = 1.500917606070914
$ws.Range("C25").Value = 0.259382267025984
$ws.Range("D25").Value = 0.521667709386449
$ws.Range("E25").Value = 0.1807368274708914
$ws.Range("G25").Value = 1.022460071671659
$ws.Range("H25").Value = 1.045379826591841
$ws.Range("J25").Value = 0.08870355048572165
$ws.Range("L25").Value = 0.3733109808283501
$ws.Range("M25").Value = 0.3707289645535141
$ws.Range("O25").Value = 4.192419313442485
